$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

$ws.Range("D2").Value = '35.125.48'
$ws.Range("E2").Value = '  -0.40%  '
$ws.Range("D3").Value = '1.813.34'
$ws.Range("E3").Value = '  -1.47%  '
$ws.Range("E4").Value = '  +0.65%  '
$ws.Range("D5").Value = '233.36'
$ws.Range("E5").Value = '  +2.31%  '
$ws.Range("D6").Value = '0.613'
$ws.Range("E6").Value = '  +0.22%  '
$ws.Range("E7").Value = '  +0.68%  '
$ws.Range("D8").Value = '40.63'
$ws.Range("E8").Value = '  -6.33%  '
$ws.Range("E9").Value = '  +5.87%  '
$ws.Range("D10").Value = '0.0685'
$ws.Range("E10").Value = '  -1.25%  '
$ws.Range("E11").Value = '  -0.47%  '
$ws.Range("E12").Value = '  -1.47%  '
$ws.Range("D13").Value = '1.822.73'
$ws.Range("E13").Value = '  -0.93%  '
$ws.Range("D14").Value = '0.664'
$ws.Range("E14").Value = '  +0.46%  '
$ws.Range("E15").Value = '  -4.83%  '
$ws.Range("E16").Value = '  -1.78%  '
$ws.Range("D17").Value = '35.098.77'
$ws.Range("E17").Value = '  -0.27%  '
$ws.Range("D18").Value = '69.63'
$ws.Range("E18").Value = '  -0.67%  '
$ws.Range("D19").Value = '0.0₃0791'
$ws.Range("E19").Value = '  -0.37%  '
$ws.Range("D20").Value = '238.78'
$ws.Range("E20").Value = '  -2.79%  '
$ws.Range("D21").Value = '11.92'
$ws.Range("E21").Value = '  -1.70%  '
$ws.Range("D22").Value = '4.71'
$ws.Range("E22").Value = '  +0.10%  '
$ws.Range("E23").Value = '  +0.70%  '
$ws.Range("D24").Value = '2.25'
$ws.Range("E24").Value = '  +2.73%  '
$ws.Range("D25").Value = '171.88'
$ws.Range("E25").Value = '  -0.26%  '
$ws.Range("D26").Value = '7.84'
$ws.Range("E26").Value = '  -1.10%  '
$ws.Range("E28").Value = '  -1.36%  '
$ws.Range("E29").Value = '  +21.84%  '
$ws.Range("E30").Value = '  +0.66%  '
$ws.Range("D31").Value = '4.18'
$ws.Range("E31").Value = '  +6.02%  '
$ws.Range("D32").Value = '3.329.34'
$ws.Range("E32").Value = '  -7.90%  '
$ws.Range("D33").Value = '0.0554'
$ws.Range("E33").Value = '  +2.98%  '
$ws.Range("D34").Value = '4.03'
$ws.Range("E34").Value = '  -0.35%  '
$ws.Range("E35").Value = '  -5.82%  '
$ws.Range("E36").Value = '  +5.13%  '
$ws.Range("D37").Value = '92.43'
$ws.Range("E37").Value = '  +2.26%  '
$ws.Range("D38").Value = '0.678'
$ws.Range("E38").Value = '  +0.47%  '
$ws.Range("D39").Value = '0.0194'
$ws.Range("E39").Value = '  +0.00%  '
$ws.Range("D40").Value = '1.313.13'
$ws.Range("E40").Value = '  -2.26%  '
$ws.Range("E41").Value = '  +1.69%  '
$ws.Range("E42").Value = '  -1.85%  '
$ws.Range("E43").Value = '  +1.00%  '
$ws.Range("D44").Value = '14.56'
$ws.Range("E44").Value = '  -2.35%  '
$ws.Range("D45").Value = '2.30'
$ws.Range("E45").Value = '  -5.99%  '
$ws.Range("E46").Value = '  -2.07%  '
$ws.Range("D47").Value = '6.30'
$ws.Range("E47").Value = '  +3.73%  '
$ws.Range("D48").Value = '0.0511'
$ws.Range("E48").Value = '  -1.61%  '
$ws.Range("E49").Value = '  -0.94%  '
$ws.Range("E50").Value = '  +0.63%  '
$ws.Range("E51").Value = '  +4.70%  '
